$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New work log entry for row 16 - match the style already used by the rows above
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A16").Value = 42899

$ws.Range("B16").Value = 0.5

$ws.Range("C16").Value = "Lepší response (logo pryč), mapa nezpůsobuje spodní scrollbar, rezervanto tlačítko otevírá konkrétní služby+nastavení reservanta, hezčí homepage slidery"

# Move the active selection to C16 (next empty row after the new entry)
$ws.Range("C16").Select() | Out-Null

$ws.Calculate() | Out-Null
